$d = $word.ActiveDocument

# Unicode helper characters used in the question text
$nu    = [char]0x03BD   # ν
$eta   = [char]0x03B7   # η
$minus = [char]0x2212   # − (minus sign)
$rsq   = [char]0x2019   # ’ (right single quotation mark)

# ---------------------------------------------------------------------
# 1) "Describe the possible states, initial state, transition function."
#    Collapses 3 runs (with proofErr gramStart/gramEnd around "transition")
#    into a single run.
# ---------------------------------------------------------------------
$t1 = "Describe the possible states, initial state, transition function."
$d.Content.Find.Execute($t1, $true, $false, $false, $false, $false, $true, 1, $false, $t1, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Checkers terminal-state paragraph: extend the sentence with more
#    detail about no-viable-moves and the 50-move draw rule.
#    The edit is applied as replacing "eft so either A wins or B wins. "
#    (the tail of the original sentence) with the new continuation, which
#    mirrors how the text was actually typed over the old tail.
# ---------------------------------------------------------------------
$checkersPara = $d.Paragraphs(5)
$ckStart = $checkersPara.Range.Start
$prefix = "The terminal state for checkers is when one of the players has no pieces l"
$oldTail = "eft so either A wins or B wins. "
$newTail = "eft so either A wins or B wins, or when the opponent has no viable moves (the other player has pinned them in). The game is a draw if neither player captures their opponent" + $rsq + "s piece in 50 total moves."
$tailStart = $ckStart + $prefix.Length
$tailRng = $d.Range($tailStart, $tailStart + $oldTail.Length)
$tailRng.Text = $newTail

# ---------------------------------------------------------------------
# 3) TTT terminal states paragraph: merge "TTT terminal states is" and
#    " when a player gets 3 " into a single run (drop the proofErr tags).
# ---------------------------------------------------------------------
$tttPara = $d.Paragraphs(7)
$tttStart = $tttPara.Range.Start
$firstOld = "TTT terminal states is"
$firstRng = $d.Range($tttStart, $tttStart + $firstOld.Length)
$firstRng.Text = "TTT terminal states is when a player gets 3 "
$dupStart = $tttStart + ("TTT terminal states is when a player gets 3 ").Length
$dupOld = " when a player gets 3 "
$dupRng = $d.Range($dupStart, $dupStart + $dupOld.Length)
$dupRng.Text = ""

# ---------------------------------------------------------------------
# 4) "Why is ν(A,s) = #{white checkers} − #{red checkers} ..." question:
#    collapse the 3 runs (proofErr wrapped "ν(") into a single run.
# ---------------------------------------------------------------------
$t4 = "Why is " + $nu + "(A,s) = #{white checkers} " + $minus + " #{red checkers} a valid heuristic function for checkers (knowing that A plays white and B plays red)?"
$d.Content.Find.Execute($t4, $true, $false, $false, $false, $false, $true, 1, $false, $t4, 2) | Out-Null

# ---------------------------------------------------------------------
# 5) "v best approximates the utility function ..." : merge "v" and the
#    rest into a single run (drop proofErr tags).
# ---------------------------------------------------------------------
$t5 = "v best approximates the utility function when the players have the same piece alignment (i.e. they have the same number of kings, and their pieces can make similar moves). "
$d.Content.Find.Execute($t5, $true, $false, $false, $false, $false, $true, 1, $false, $t5, 2) | Out-Null

# ---------------------------------------------------------------------
# 6) "Can you provide an example of a state s where v(A,s)>0 ..." :
#    collapse the 3 runs (proofErr wrapped "v(") into a single run.
# ---------------------------------------------------------------------
$t6 = "Can you provide an example of a state s where v(A,s)>0 and B wins in the following turn? (Hint: recall the rules for jumping in checkers)"
$d.Content.Find.Execute($t6, $true, $false, $false, $false, $false, $true, 1, $false, $t6, 2) | Out-Null

# ---------------------------------------------------------------------
# 7) Move the "_GoBack" bookmark from the end of the document to the end
#    of the (edited) checkers terminal-state paragraph.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$checkersPara2 = $d.Paragraphs(5)
$newEnd = $checkersPara2.Range.End - 2   # one char before the trailing period's end, just inside the paragraph text
$bmRange = $d.Range($newEnd, $newEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
